$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D").Insert()

# Copy number formats (and other formatting) from the now-shifted column E into the new column D
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new column D with FY2018 data
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 5485100
$ws.Range("D9").Value = 3515600
$ws.Range("D10").Value = 1969500
$ws.Range("D12").Value = 50300
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 98000
$ws.Range("D15").Value = 36100
$ws.Range("D17").Value = 4889900
$ws.Range("D18").Value = 595200
$ws.Range("D20").Value = 16300
$ws.Range("D21").Value = 761100
$ws.Range("D22").Value = 74500
$ws.Range("D23").Value = 537000
$ws.Range("D24").Value = 141500
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 395500
$ws.Range("D27").Value = 395300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -5700
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -16300
$ws.Range("D33").Value = 389600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 389600
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 262900
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 571700
$ws.Range("D44").Value = 678900
$ws.Range("D45").Value = 172600
$ws.Range("D46").Value = 1686100
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 813400
$ws.Range("D49").Value = 3327100
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 138000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 5964600
$ws.Range("D57").Value = 459000
$ws.Range("D58").Value = 525000
$ws.Range("D59").Value = 508100
$ws.Range("D60").Value = 1492100
$ws.Range("D61").Value = 1809000
$ws.Range("D62").Value = 483500
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 3786400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 1448100
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 2178200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 389600
$ws.Range("D83").Value = 149600
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 604000
$ws.Range("D91").Value = -150100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -634300
$ws.Range("D96").Value = -115200
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -6800
$ws.Range("D101").Value = -15200
$ws.Range("D102").Value = -52300
